{"js": "// Replace the date line and each \"a\u00f7b=c, r\" answer cell with its updated value.\n// Every old string below occurs exactly once in the document, so a scoped\n// search-and-replace (matchCase, no wildcards) is safe and keeps each run's\n// original formatting (rFonts/sz) untouched.\nconst replacements = [\n  [\"2025-08-25 Monday\", \"2025-08-26 Tuesday\"],\n  [\"703\u00f77=100, 3\", \"170\u00f75=34, 0\"],\n  [\"803\u00f72=401, 1\", \"160\u00f73=53, 1\"],\n  [\"785\u00f72=392, 1\", \"934\u00f79=103, 7\"],\n  [\"873\u00f73=291, 0\", \"963\u00f78=120, 3\"],\n  [\"994\u00f73=331, 1\", \"615\u00f75=123, 0\"],\n  [\"105\u00f79=11, 6\", \"156\u00f75=31, 1\"],\n  [\"357\u00f72=178, 1\", \"430\u00f78=53, 6\"],\n  [\"180\u00f78=22, 4\", \"132\u00f74=33, 0\"],\n  [\"840\u00f76=140, 0\", \"478\u00f79=53, 1\"],\n  [\"776\u00f72=388, 0\", \"689\u00f79=76, 5\"],\n  [\"827\u00f73=275, 2\", \"481\u00f76=80, 1\"],\n  [\"495\u00f75=99, 0\", \"170\u00f79=18, 8\"],\n  [\"459\u00f77=65, 4\", \"752\u00f79=83, 5\"],\n  [\"629\u00f75=125, 4\", \"280\u00f74=70, 0\"],\n  [\"246\u00f74=61, 2\", \"707\u00f73=235, 2\"],\n  [\"691\u00f76=115, 1\", \"632\u00f74=158, 0\"],\n  [\"233\u00f74=58, 1\", \"350\u00f77=50, 0\"],\n  [\"714\u00f78=89, 2\", \"585\u00f78=73, 1\"],\n  [\"549\u00f75=109, 4\", \"399\u00f73=133, 0\"],\n  [\"289\u00f78=36, 1\", \"381\u00f75=76, 1\"],\n  [\"732\u00f78=91, 4\", \"943\u00f79=104, 7\"],\n  [\"344\u00f72=172, 0\", \"228\u00f78=28, 4\"],\n  [\"783\u00f77=111, 6\", \"333\u00f78=41, 5\"],\n  [\"864\u00f73=288, 0\", \"230\u00f79=25, 5\"],\n  [\"217\u00f75=43, 2\", \"949\u00f78=118, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"a\u00f7b=c, r\" answer cell with its updated\n# value. Every old string occurs exactly once in the document, so a\n# Find/Replace-one pass per pair is unambiguous and preserves each run's\n# original formatting (rFonts/sz) since Find.Execute only swaps the text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-08-25 Monday\", \"2025-08-26 Tuesday\"),\n    @(\"703\u00f77=100, 3\", \"170\u00f75=34, 0\"),\n    @(\"803\u00f72=401, 1\", \"160\u00f73=53, 1\"),\n    @(\"785\u00f72=392, 1\", \"934\u00f79=103, 7\"),\n    @(\"873\u00f73=291, 0\", \"963\u00f78=120, 3\"),\n    @(\"994\u00f73=331, 1\", \"615\u00f75=123, 0\"),\n    @(\"105\u00f79=11, 6\", \"156\u00f75=31, 1\"),\n    @(\"357\u00f72=178, 1\", \"430\u00f78=53, 6\"),\n    @(\"180\u00f78=22, 4\", \"132\u00f74=33, 0\"),\n    @(\"840\u00f76=140, 0\", \"478\u00f79=53, 1\"),\n    @(\"776\u00f72=388, 0\", \"689\u00f79=76, 5\"),\n    @(\"827\u00f73=275, 2\", \"481\u00f76=80, 1\"),\n    @(\"495\u00f75=99, 0\", \"170\u00f79=18, 8\"),\n    @(\"459\u00f77=65, 4\", \"752\u00f79=83, 5\"),\n    @(\"629\u00f75=125, 4\", \"280\u00f74=70, 0\"),\n    @(\"246\u00f74=61, 2\", \"707\u00f73=235, 2\"),\n    @(\"691\u00f76=115, 1\", \"632\u00f74=158, 0\"),\n    @(\"233\u00f74=58, 1\", \"350\u00f77=50, 0\"),\n    @(\"714\u00f78=89, 2\", \"585\u00f78=73, 1\"),\n    @(\"549\u00f75=109, 4\", \"399\u00f73=133, 0\"),\n    @(\"289\u00f78=36, 1\", \"381\u00f75=76, 1\"),\n    @(\"732\u00f78=91, 4\", \"943\u00f79=104, 7\"),\n    @(\"344\u00f72=172, 0\", \"228\u00f78=28, 4\"),\n    @(\"783\u00f77=111, 6\", \"333\u00f78=41, 5\"),\n    @(\"864\u00f73=288, 0\", \"230\u00f79=25, 5\"),\n    @(\"217\u00f75=43, 2\", \"949\u00f78=118, 5\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue=1, wdReplaceOne=2\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
